$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells we touch remain text (matches original inlineStr type)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "40.914.21"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -2.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.388.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.63%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.67"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.69"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.92%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.71%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -4.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0841"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "31.28"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.07%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.756.38"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.65%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.42%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.392.93"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.65%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.770"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.66%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "40.868.44"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0916"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -3.76%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -4.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "69.35"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.94"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.68%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "233.44"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.65"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.97%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -6.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.03"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.98%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.64%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.39"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -4.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.15"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -7.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "153.85"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -5.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0737"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.81%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.01%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.82"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.10"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -8.06%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.80%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -7.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.85"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.04%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.974.96"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.60"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -7.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.63"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.35%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.76"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.624.55"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "94.09"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.48"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.07%  "
